# Swap the order of the two slides in the "exec summary" portion of the
# deck: the slide titled "Frame the Problem" (currently 4th) should come
# right before the slide titled "Terms Defined" (currently 3rd).
#
# Moving slide 4 up to slide position 3 shifts the former slide 3 down to
# position 4, i.e. a simple adjacent swap of slides 3 and 4 - everything
# that lives on those slides (shapes, text, speaker notes, etc.) travels
# together with them.

$p = $ppt.ActivePresentation

$frameTheProblem = $p.Slides.Item(4)
$frameTheProblem.MoveTo(3)
